$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Activate()

# The cell E8 held the "Good Morning" greeting text; update it to "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Leave E8 selected, matching the saved selection state in the target file.
$ws.Range("E8").Select()
